$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.280.96"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.524.68"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +4.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0821"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "2.917.12"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "2.524.39"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.862"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "48.168.66"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "0.0₃0946"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.146"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0793"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0300"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "2.007.53"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("E47").Value = "  +6.48%  "
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.01%  "
